# Apply updated cryptos table values (price + 1h volume%) for rows 2-51.
# Rows 33/34 additionally swap Coin/Link (Hedera <-> EthereumClassic reordered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.018.66"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "2.269.14"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'318.41"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'102.29"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.571"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "'38.71"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'7.85"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "2.613.95"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'0.875"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").Value = "2.265.42"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "44.005.79"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "'14.43"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("D21").Value = "'6.67"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "'66.03"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").Value = "'3.21"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "'238.87"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").Value = "'10.24"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'39.26"
$ws.Range("E29").Value = "  +15.63%  "
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").Value = "'6.51"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").Value = "'163.71"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.46"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0880"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "'2.71"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("D36").Value = "'3.25"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'2.04"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'4.56"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("D42").Value = "'15.76"
$ws.Range("E42").Value = "  +29.18%  "
$ws.Range("D43").Value = "'0.0326"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "1.775.48"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "'85.05"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").Value = "'5.40"
$ws.Range("E48").Value = "  -2.21%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("E49").Value = "  +4.08%  "
$ws.Range("D50").Value = "'59.65"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("D51").Value = "'74.99"
$ws.Range("E51").Value = "  -4.02%  "
